$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
# "Volume 32   Number  8" -> "...Number  9" (replace the last run "8" with "9")
$volChars = $ws.Range("A8").Characters(21, 1)
$volChars.Text = "9"

# "Report Covering the Week  2/17/2025  Through  2/23/2025"
#   -> "...2/24/2025  Through  3/2/2025"
$weekStart = $ws.Range("C9").Characters(27, 9)
$weekStart.Text = "2/24/2025"
$weekEnd = $ws.Range("C9").Characters(47, 9)
$weekEnd.Text = "3/2/2025"

# --- Weekly crime statistics table (rows 14-30) ---
# Row 14
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("G14").Value = 2
$ws.Range("J14").Value = 5

# Row 15
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 4
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -60
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -69.230769230769

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -36.363636363636
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = -31.372549019607
$ws.Range("L16").Value = -20.454545454545
$ws.Range("M16").Value = -39.655172413793
$ws.Range("N16").Value = -90.74074074074

# Row 17
$ws.Range("C17").Value = 22
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = 29.411764705882
$ws.Range("F17").Value = 59
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = 31.111111111111
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = 8.910891089108
$ws.Range("L17").Value = 0.91743119266
$ws.Range("M17").Value = 44.736842105263
$ws.Range("N17").Value = -30.379746835443

# Row 18
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = 17
$ws.Range("J18").Value = 19
$ws.Range("K18").Value = -10.526315789473
$ws.Range("L18").Value = -43.333333333333
$ws.Range("M18").Value = -61.363636363636
$ws.Range("N18").Value = -86.71875

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 10.714285714285
$ws.Range("I19").Value = 50
$ws.Range("J19").Value = 67
$ws.Range("K19").Value = -25.373134328358
$ws.Range("L19").Value = -23.076923076923
$ws.Range("M19").Value = -13.793103448275
$ws.Range("N19").Value = -55.75221238938

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -41.666666666666
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = -44.827586206896
$ws.Range("L20").Value = -11.111111111111
$ws.Range("M20").Value = -15.78947368421
$ws.Range("N20").Value = -85.185185185185

# Row 21
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = 5.128205128205
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = -0.813008130081
$ws.Range("I21").Value = 232
$ws.Range("J21").Value = 276
$ws.Range("K21").Value = -15.942028985507
$ws.Range("L21").Value = -16.546762589928
$ws.Range("M21").Value = -12.121212121212
$ws.Range("N21").Value = -74.421168687982

# Row 22
$ws.Range("D22").Value = 4
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = -71.428571428571
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = -14.285714285714

# Row 23
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = 12.5
$ws.Range("F23").Value = 21
$ws.Range("G23").Value = 33
$ws.Range("H23").Value = -36.363636363636
$ws.Range("I23").Value = 47
$ws.Range("J23").Value = 61
$ws.Range("K23").Value = -22.950819672131
$ws.Range("L23").Value = -29.850746268656
$ws.Range("M23").Value = 46.875

# Row 24
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -17.391304347826
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = -35.514018691588
$ws.Range("I24").Value = 148
$ws.Range("J24").Value = 211
$ws.Range("K24").Value = -29.857819905213
$ws.Range("L24").Value = -36.206896551724
$ws.Range("M24").Value = 0

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = -40.625
$ws.Range("I25").Value = 34
$ws.Range("J25").Value = 62
$ws.Range("K25").Value = -45.16129032258
$ws.Range("L25").Value = -43.333333333333

# Row 26
$ws.Range("D26").Value = 23
$ws.Range("E26").Value = -39.130434782608
$ws.Range("F26").Value = 64
$ws.Range("G26").Value = 63
$ws.Range("H26").Value = 1.587301587301
$ws.Range("I26").Value = 141
$ws.Range("J26").Value = 128
$ws.Range("K26").Value = 10.15625
$ws.Range("L26").Value = -12.422360248447
$ws.Range("M26").Value = -25.396825396825

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 4
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -33.333333333333
$ws.Range("L27").Value = -63.636363636363

# Row 28
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = 45.454545454545
$ws.Range("L28").Value = -15.78947368421

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 4
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = -60
$ws.Range("L29").Value = -33.333333333333
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -88.888888888888

# Row 30
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 50
$ws.Range("I30").Value = 4
$ws.Range("J30").Value = 9
$ws.Range("K30").Value = -55.555555555555
$ws.Range("L30").Value = -33.333333333333
$ws.Range("M30").Value = -33.333333333333
$ws.Range("N30").Value = -88.888888888888
